# Rename the existing sheet "Sheet2" -> "Matrix", then add a new
# worksheet "Workers" right after it, filling column A with the
# list of agent names (matching the style used for the data rows
# on the Matrix sheet).

$wb = $excel.ActiveWorkbook

$matrix = $wb.Worksheets.Item(1)
$matrix.Name = "Matrix"

$workers = $wb.Worksheets.Add($null, $matrix)
$workers.Name = "Workers"

$agents = @(
    "Agent",
    "Silva, Flavio",
    "Gomez, Manuel",
    "Bertone, Ignacio",
    "Pereira,Eugenia",
    "Taborda, Lucia",
    "Gallinar,Romina",
    "Oliver, Matilde",
    "Alvez, Eugenia",
    "Florin, Steban",
    "Rodriguez, Ginni",
    "Gill, Angela"
)

for ($i = 0; $i -lt $agents.Length; $i++) {
    $workers.Cells.Item($i + 1, 1).Value = $agents[$i]
}

# Match the formatting used for data cells on the Matrix sheet (e.g. B4)
# - thin border all around, Arial 10 font - by copying the format over.
$matrix.Range("B4").Copy()
$workers.Range("A1:A12").PasteSpecial(-4122)

# Widen column A to fit the agent names, matching the Workers sheet layout.
$workers.Columns.Item(1).ColumnWidth = 12.5

$workers.Range("A1").Select()
